# Split the old "main" sheet (9 rows mixing EU/France/UK passenger time
# values) into three dedicated per-country sheets - eu25, france, uk - and
# leave the "additional" sheet's data untouched. New tab order:
# eu25, france, uk, additional.

$wb = $excel.ActiveWorkbook

$mainSheet = $wb.Worksheets.Item("main")
$additionalSheet = $wb.Worksheets.Item("additional")

# Drop the old combined sheet - its data is being redistributed below.
$mainSheet.Delete()

# Re-create the three new sheets, inserted right before "additional" so the
# final left-to-right order is eu25, france, uk, additional.
$eu25 = $wb.Worksheets.Add($additionalSheet)
$eu25.Name = "eu25"

$france = $wb.Worksheets.Add($null, $eu25)
$france.Name = "france"

$uk = $wb.Worksheets.Add($null, $france)
$uk.Name = "uk"

# --- eu25: Personal / Business -------------------------------------------
$eu25.Range("A1").Value = "purpose"
$eu25.Range("B1").Value = "value"
$eu25.Range("A2").Value = "Personal"
$eu25.Range("B2").Value = 20.47
$eu25.Range("A3").Value = "Business"
$eu25.Range("B3").Value = 49.98

# --- france: Personal - holiday / Personal - other / Business / All purpose
$france.Range("A1").Value = "purpose"
$france.Range("B1").Value = "value"
$france.Range("A2").Value = "Personal - holiday"
$france.Range("B2").Value = 59.8
$france.Range("A3").Value = "Personal - other"
$france.Range("B3").Value = 61.2
$france.Range("A4").Value = "Business"
$france.Range("B4").Value = 83.5
$france.Range("A5").Value = "All purpose"
$france.Range("B5").Value = 62.1

# --- uk: Leisure / UK business / Foreign business -------------------------
$uk.Range("A1").Value = "purpose"
$uk.Range("B1").Value = "value"
$uk.Range("A2").Value = "Leisure"
$uk.Range("B2").Value = 8.6
$uk.Range("A3").Value = "UK business"
$uk.Range("B3").Value = 63.8
$uk.Range("A4").Value = "Foreign business"
$uk.Range("B4").Value = 60.7

# Selections mirroring the post-edit session state on each sheet.
$eu25.Range("K12").Select()
$france.Range("B5").Select()
$uk.Range("A1:B1").Select()

# Move "additional" to the end so the tab order is eu25, france, uk, additional.
$additionalSheet = $wb.Worksheets.Item("additional")
$ukSheet = $wb.Worksheets.Item("uk")
$additionalSheet.Move($null, $ukSheet)

# eu25 is the active/selected tab.
$eu25sheet = $wb.Worksheets.Item("eu25")
$eu25sheet.Activate()
$eu25sheet.Range("K12").Select()
